# Add benchmark results for two more mobile Ryzen 5 laptop CPUs
# (Ryzen 5 4500U / ThinkPad L14, and Ryzen 5 PRO 4650U / ThinkPad T14s)
# into the sorted CPU benchmark table on the "List1" worksheet.
#
# The table is kept sorted by column O (rating, ascending), so the two
# new rows are inserted in-place at their correctly sorted positions
# (row 10 and row 13), pushing the existing rows below them down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at their sorted positions. Inserting at row 10
# first, then at row 13 (post-insert numbering), reproduces the final
# layout: old row 10 -> 11, old row 11 -> 12, new row -> 13, old row 12 -> 14, ...
$ws.Rows("10:10").Insert()
$ws.Rows("13:13").Insert()

# Fill row 13 (Ryzen 5 PRO 4650U / ThinkPad T14s) before row 10 so new
# shared strings are appended in the same order as the source edit.
$ws.Cells.Item(13, 1).Value = "AMD"
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = "Ryzen 5 PRO 4650U"
$ws.Cells.Item(13, 4).Value = 25
$ws.Cells.Item(13, 5).Value = 6
$ws.Cells.Item(13, 6).Value = 2
$ws.Cells.Item(13, 7).Value = 2.1
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(13, 9).Value = "x86-64"
$ws.Cells.Item(13, 11).Value = 16
$ws.Cells.Item(13, 12).Value = 2
$ws.Cells.Item(13, 13).Value = "DDR4"
$ws.Cells.Item(13, 14).Value = 3200
$ws.Cells.Item(13, 15).Value = 0.51
$ws.Cells.Item(13, 16).Value = 1.05
$ws.Cells.Item(13, 17).Value = 2.2
$ws.Cells.Item(13, 18).Value = 4.32
$ws.Cells.Item(13, 19).Value = "ThinkPad T14s"

# Fill row 10 (Ryzen 5 4500U / ThinkPad L14)
$ws.Cells.Item(10, 1).Value = "AMD"
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "Ryzen 5 4500U"
$ws.Cells.Item(10, 4).Value = 15
$ws.Cells.Item(10, 5).Value = 6
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.3
$ws.Cells.Item(10, 8).Value = 4
$ws.Cells.Item(10, 9).Value = "x86-64"
$ws.Cells.Item(10, 11).Value = 32
$ws.Cells.Item(10, 12).Value = 2
$ws.Cells.Item(10, 13).Value = "DDR4"
$ws.Cells.Item(10, 14).Value = 3200
$ws.Cells.Item(10, 15).Value = 0.46
$ws.Cells.Item(10, 16).Value = 0.91
$ws.Cells.Item(10, 17).Value = 1.88
$ws.Cells.Item(10, 18).Value = 3.75
$ws.Cells.Item(10, 19).Value = "ThinkPad L14"

# Match the author's final selection (cell J49) left after the edit.
$ws.Range("J49").Select() | Out-Null
